$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.378.94"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.047.44"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.01%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.87"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "619.00"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.96%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +0.34%  "

$ws.Range("E9").Value = "  +6.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.044.72"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.22%  "

$ws.Range("E11").Value = "  -1.73%  "

$ws.Range("E12").Value = "  -0.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.22"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.604.44"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.94"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.67%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.265.09"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000193"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.028.22"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.56"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.93"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.28"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.44"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.05%  "

$ws.Range("E23").Value = "  +2.00%  "

$ws.Range("E24").Value = "  +2.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.42"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  +1.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.85"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.84%  "

$ws.Range("E29").Value = "  +1.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.996"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.44%  "

$ws.Range("E31").Value = "  +1.09%  "

$ws.Range("E32").Value = "  +1.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "496.67"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.39%  "

$ws.Range("E34").Value = "  +5.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.125"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.66"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.48"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.58%  "

$ws.Range("E39").Value = "  +1.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "192.14"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.378"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.14%  "

$ws.Range("E42").Value = "  -9.75%  "

$ws.Range("E44").Value = "  +4.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.784"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +19.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.25"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.97%  "

$ws.Range("E47").Value = "  +2.65%  "

$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.44"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.595"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.87"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.44%  "
